$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$result = $find.Execute("the place I live is not ", $true, $false, $false, $false, $false, $true, 1, $false, "where I live is, at minimum, 7 kilometers away from any point of interest. ", 2)
Write-Output "Replace result: $result"
